$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) for rows 2 through 41: 45703 -> 45704
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 3).Value = 45704
}

# Swap row 40 and row 41 values for column A (Beteckning) and column G (Area (ha))
$ws.Cells.Item(40, 1).Value = "A 3005-2025"
$ws.Cells.Item(41, 1).Value = "A 3002-2025"

$ws.Cells.Item(40, 7).Value = 2.2
$ws.Cells.Item(41, 7).Value = 2.4
